# CSE111_ProjectUpdated.pptx - "Add files via upload" edit replay
#
# 1) Shrink the "ra_*" lookup-table (slide 4, table shape id=100) from a
#    646125 EMU wide column to 544800 EMU, and rename the rating columns
#    to the "ad_*" (Adult) columns.
# 2) Add a new "repo link" text box under the Demo title on the last
#    slide (slide 7).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 4 - shrink + relabel the small lookup table (shape id 100)
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)

$tableShape = $null
foreach ($shp in $slide4.Shapes) {
    if ($shp.Id -eq 100) {
        $tableShape = $shp
    }
}

$tbl = $tableShape.Table

# 646125 EMU -> 544800 EMU (EMU / 12700 = points)
$tbl.Columns.Item(1).Width = 544800 / 12700

$tbl.Cell(1, 1).Shape.TextFrame.TextRange.Text = "Adult"
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "ad_id"
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text = "ad_ratings"

# ---------------------------------------------------------------------
# 2. Slide 7 - add the GitHub repo link text box under the "Demo" title
# ---------------------------------------------------------------------
$slide7 = $p.Slides.Item($p.Slides.Count)

# Burn shape-id 3 on a throwaway textbox so the real one lands on id 4 /
# "TextBox 3", matching what PowerPoint assigned when it created this
# shape (ids keep climbing even across deletes).
$placeholder = $slide7.Shapes.AddTextbox(1, 0, 0, 10, 10)
$placeholder.Delete()

$linkBox = $slide7.Shapes.AddTextbox(1, 2538430 / 12700, 3570135 / 12700, 4067139 / 12700, 307777 / 12700)
$linkBox.Fill.Visible = $false

$linkText = $linkBox.TextFrame.TextRange
$linkText.Text = "https://"
$linkText.InsertAfter("github.com")
$linkText.InsertAfter("/")
$linkText.InsertAfter("winstonlou")
$linkText.InsertAfter("/CSE111-Project.git")

# Single-line, auto-fit-to-text box (no wrapping), like a typed hyperlink.
$linkBox.TextFrame.WordWrap = $false
$linkBox.TextFrame.AutoSize = 1

# AutoSize recomputes Height off an 18pt fallback instead of this deck's
# real 14pt default text size, so pin it back to the authored 307777 EMU
# (add half an EMU-in-points so the engine's internal float32 round-trip
# floors back to the right integer EMU instead of the one below it).
$linkBox.Height = (307777 + 0.5) / 12700
